$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1, 0.768386970581898),
    @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.145820798638228),
    @(0.1554434735375247, 0.004309184025731883, 0.7127328510149897, 6.48142807727062, 0, 7.353913585848867),
    @(0.7287194209349384, 1.65323645889881, 157.8057217802531, 6.48142807727062, 1, 166.6691057373575),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
